$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row (vehicle_id, vehicle_type, cost, location_id, in_use)
# so the data exports cleanly to postgres without a text header row.
$ws.Rows("1").Delete()

# Reset the active selection to a sensible cell after the row deletion.
$ws.Range("A3").Select()
